$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow editing, then restore protection at the end.
$ws.Unprotect("D382")

# Update the confidential/date disclosure text in cell A42 (shared string), bumping the
# "as of" date from 2021-05-24 to 2021-05-25.
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

# Update the Weight (column D) and Percent Change (column E) values for each holdings row.
$ws.Range("D2").Value = 0.05741310638279992
$ws.Range("E2").Value = -0.001573564122738014
$ws.Range("D3").Value = 0.05268899719568391
$ws.Range("E3").Value = 0.003748305287502918
$ws.Range("D4").Value = 0.3164218208252642
$ws.Range("D5").Value = 0.03408869706715494
$ws.Range("E5").Value = 0.004332833074986553
$ws.Range("D6").Value = 0.03131715664723763
$ws.Range("E6").Value = -0.001308215593929796
$ws.Range("D7").Value = 0.0309238485574764
$ws.Range("E7").Value = -0.01033386327503982
$ws.Range("D8").Value = 0.02866610884990475
$ws.Range("E8").Value = -0.002755790090882448
$ws.Range("D9").Value = 0.02382707470279975
$ws.Range("E9").Value = 0.004091422121896171
$ws.Range("D10").Value = 0.02480278131009202
$ws.Range("E10").Value = 0.0007750821671805852
$ws.Range("D11").Value = 0.02387172105380096
$ws.Range("E11").Value = 0.00973415888858109
$ws.Range("D12").Value = 0.02328711648127927
$ws.Range("E12").Value = -0.01454374853389639
$ws.Range("D13").Value = 0.0201424376172295
$ws.Range("E13").Value = 0.01067064425448905
$ws.Range("D14").Value = 0.02173909617927033
$ws.Range("E14").Value = 0.003406784575239241
$ws.Range("D15").Value = 0.01984441009301674
$ws.Range("E15").Value = 0.006066573497649719
$ws.Range("D16").Value = 0.02158519758111323
$ws.Range("E16").Value = 0.005110110719065641
$ws.Range("D17").Value = 0.01955930374803491
$ws.Range("E17").Value = -0.01385681293302543
$ws.Range("D18").Value = 0.01436078768769195
$ws.Range("E18").Value = -0.001580056179775302
$ws.Range("D19").Value = 0.01637186943726164
$ws.Range("E19").Value = 0.0007186489399928497
$ws.Range("D20").Value = 0.01496440635322827
$ws.Range("E20").Value = -0.008775008775008808
$ws.Range("D21").Value = 0.01628131613240742
$ws.Range("E21").Value = -0.02264720684448918
$ws.Range("D22").Value = 0.01274133322408109
$ws.Range("E22").Value = -0.002885693555834012
$ws.Range("D23").Value = 0.01496755786035777
$ws.Range("E23").Value = -0.0001824817518247812
$ws.Range("D24").Value = 0.01343855165136348
$ws.Range("E24").Value = -0.007731092436974896
$ws.Range("D25").Value = 0.01402158047032042
$ws.Range("E25").Value = -0.001498407941562085
$ws.Range("D26").Value = 0.01355347661135247
$ws.Range("E26").Value = 0.009463722397476282
$ws.Range("D27").Value = 0.01267473137341105
$ws.Range("E27").Value = 0.001939430105922568
$ws.Range("D28").Value = 0.01331196611499534
$ws.Range("E28").Value = -0.01638257575757562
$ws.Range("D29").Value = 0.01421896986686458
$ws.Range("E29").Value = -0.01205727204220042
$ws.Range("D30").Value = 0.01334222058343852
$ws.Range("E30").Value = -0.006928697404887885
$ws.Range("D31").Value = 0.01242534210923019
$ws.Range("E31").Value = 0.00304362529590807
$ws.Range("D32").Value = 0.01339401035059992
$ws.Range("E32").Value = -0.01164696747476479
$ws.Range("D33").Value = 0.01245318042220741
$ws.Range("E33").Value = -0.01935973680880687
$ws.Range("D34").Value = 0.006560177240760964
$ws.Range("E34").Value = 0.002289905201127151
$ws.Range("D35").Value = 0.005282976451413478
$ws.Range("E35").Value = -0.003102008351560936
$ws.Range("D36").Value = 0.005403153923284964
$ws.Range("E36").Value = 0.005754948088812606
$ws.Range("D37").Value = 0.005229715980924979
$ws.Range("E37").Value = 0.014563204306691
$ws.Range("D38").Value = 0.004823801862645765
$ws.Range("E38").Value = 0.000762211720638506
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 0.001835122601506312

# Restore sheet protection to match original state.
$ws.Protect("D382", $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
